$d = $word.ActiveDocument

# Locate the paragraph right after "LOM3101: Mecânica dos Materiais (Requisito)"
# and delete the following three paragraphs:
#   1) the blank paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
# The blank paragraph that follows those (right before the page-break paragraph) stays intact.

$count = $d.Paragraphs.Count
$startIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq "LOM3101: Mecânica dos Materiais (Requisito)") {
        $startIndex = $i
        break
    }
}

if ($startIndex -eq -1) {
    throw "Could not find the 'LOM3101' paragraph"
}

# The three paragraphs to remove are the ones immediately following $startIndex.
$p1 = $d.Paragraphs.Item($startIndex + 1)
$p3 = $d.Paragraphs.Item($startIndex + 3)

$rangeStart = $p1.Range.Start
$rangeEnd = $p3.Range.End

$r = $d.Range($rangeStart, $rangeEnd)
$r.Delete()

Write-Output "Deleted paragraphs after index $startIndex"
